$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.6720632314682007
$ws.Range("B1").Value = 0.5633464455604553
$ws.Range("C1").Value = 0.4193816781044006
$ws.Range("D1").Value = 0.3945804536342621
$ws.Range("E1").Value = 0.4147497415542603
